$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 27-32: a continuation of the existing 3-row repeating cycle
# (same data as rows 3-5), appended twice more.
$newRows = @(
    @("AA202400050262247", "050", "ANWA ODUMO", "98765422789", "ANWA ODUMO", "Canada", "Personal Travel Allowance (PTA)", "REPAYMENT OF FOREIGN LOANS", "2025-02-23", "None", "100.0", "0"),
    @("AA202400050262134", "050", "ANWA ODUMO", "98765422789", "A and T Medical services", "United States of America", "Maintenance/Upkeep", "REPAYMENT OF FOREIGN LOANS", "2025-02-21", "USD", "10.0", "0"),
    @("AA202400050262413", "050", "ANWA ODUMO", "98765422789", "A and T Medical services", "United States of America", "CONSULTANCY FEES", "REPAYMENT OF FOREIGN LOANS", "2025-02-09", "USD", "100.0", "0"),
    @("AA202400050262247", "050", "ANWA ODUMO", "98765422789", "ANWA ODUMO", "Canada", "Personal Travel Allowance (PTA)", "REPAYMENT OF FOREIGN LOANS", "2025-02-23", "None", "100.0", "0"),
    @("AA202400050262134", "050", "ANWA ODUMO", "98765422789", "A and T Medical services", "United States of America", "Maintenance/Upkeep", "REPAYMENT OF FOREIGN LOANS", "2025-02-21", "USD", "10.0", "0"),
    @("AA202400050262413", "050", "ANWA ODUMO", "98765422789", "A and T Medical services", "United States of America", "CONSULTANCY FEES", "REPAYMENT OF FOREIGN LOANS", "2025-02-09", "USD", "100.0", "0")
)

# Columns that look numeric/date-like ("050", "98765422789", "2025-02-23",
# "100.0", "0") must be forced to plain text so Excel's automatic type
# inference doesn't turn them into numbers / date serials. Column A (form
# numbers like "AA202400050262247") and the other free-text columns are
# already safe because they don't parse as numbers.
$textColumns = @(2, 4, 9, 11, 12)

$startRow = 27
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowValues = $newRows[$i]
    $r = $startRow + $i

    foreach ($c in $textColumns) {
        $ws.Cells.Item($r, $c).NumberFormat = "@"
    }

    for ($c = 1; $c -le $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }

    # Drop the temporary "@" text format again so the cell keeps the
    # workbook's default (unstyled) look, matching the other data rows.
    foreach ($c in $textColumns) {
        $ws.Cells.Item($r, $c).ClearFormats()
    }
}
